$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"
$ws.Range("B4").Value = "Pabellón De Arteaga"
$ws.Range("B5").Value = "Rincón De Romos"
$ws.Range("B21").Value = "Amatenango De La Frontera"
$ws.Range("B24").Value = "Bejucal De Ocampo"
$ws.Range("B31").Value = "Comitán De Domínguez"
$ws.Range("B47").Value = "Mazapa De Madero"
$ws.Range("B55").Value = "Salto De Agua"
$ws.Range("B56").Value = "San Cristóbal De Las Casas"
$ws.Range("B80").Value = "Hidalgo Del Parral"
$ws.Range("D97").Value = 0.0009302325581395348
$ws.Range("A103").Value = "Ciudad De México"
$ws.Range("B107").Value = "Cuajimalpa De Morelos"
$ws.Range("B121").Value = "Coneto De Comonfort"
$ws.Range("A130").Value = "Estado De México"
$ws.Range("B130").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B133").Value = "Almoloya De Alquisiras"
$ws.Range("B134").Value = "Almoloya De Juárez"
$ws.Range("B135").Value = "Almoloya Del Río"
$ws.Range("B142").Value = "Atizapán De Zaragoza"
$ws.Range("B149").Value = "Chapa De Mota"
$ws.Range("B152").Value = "Coacalco De Berriozábal"
$ws.Range("B159").Value = "Ecatepec De Morelos"
$ws.Range("B164").Value = "Ixtapan De La Sal"
$ws.Range("B165").Value = "Ixtapan Del Oro"
$ws.Range("D172").Value = 0.0009302325581395348
$ws.Range("B177").Value = "Naucalpan De Juárez"
$ws.Range("D179").Value = 0.009418604651162793
$ws.Range("B189").Value = "San Antonio La Isla"
$ws.Range("B190").Value = "San Felipe Del Progreso"
$ws.Range("B191").Value = "San Simón De Guerrero"
$ws.Range("B201").Value = "Tenango Del Valle"
$ws.Range("B212").Value = "Tlalnepantla De Baz"
$ws.Range("B216").Value = "Valle De Bravo"
$ws.Range("B217").Value = "Valle De Chalco Solidaridad"
$ws.Range("B218").Value = "Villa De Allende"
$ws.Range("B228").Value = "San Miguel De Allende"
$ws.Range("B229").Value = "Apaseo El Alto"
$ws.Range("B230").Value = "Apaseo El Grande"
$ws.Range("B235").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B247").Value = "San Diego De La Unión"
$ws.Range("B249").Value = "San Francisco Del Rincón"
$ws.Range("B251").Value = "San Luis De La Paz"
$ws.Range("B252").Value = "Silao De La Victoria"
$ws.Range("B257").Value = "Acapulco De Juárez"
$ws.Range("B260").Value = "Ajuchitlán Del Progreso"
$ws.Range("B261").Value = "Alcozauca De Guerrero"
$ws.Range("B265").Value = "Atenango Del Río"
$ws.Range("B266").Value = "Atlamajalcingo Del Monte"
$ws.Range("B268").Value = "Atoyac De Álvarez"
$ws.Range("B269").Value = "Ayutla De Los Libres"
$ws.Range("B271").Value = "Chilapa De Álvarez"
$ws.Range("B272").Value = "Chilpancingo De Los Bravo"
$ws.Range("B273").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B278").Value = "Coyuca De Benítez"
$ws.Range("B279").Value = "Coyuca De Catalán"
$ws.Range("B283").Value = "Cuetzala Del Progreso"
$ws.Range("B284").Value = "Cutzamala De Pinzón"
$ws.Range("B290").Value = "Huitzuco De Los Figueroa"
$ws.Range("B291").Value = "Iguala De La Independencia"
$ws.Range("B293").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B294").Value = "Zihuatanejo De Azueta"
$ws.Range("B296").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B299").Value = "Mártir De Cuilapan"
$ws.Range("B311").Value = "Taxco De Alarcón"
$ws.Range("B313").Value = "Técpan De Galeana"
$ws.Range("B315").Value = "Tepecoacuilco De Trujano"
$ws.Range("B317").Value = "Tixtla De Guerrero"
$ws.Range("B321").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B322").Value = "Tlapa De Comonfort"
$ws.Range("B334").Value = "Agua Blanca De Iturbide"
$ws.Range("B341").Value = "Atotonilco El Grande"
$ws.Range("B347").Value = "Cuautepec De Hinojosa"
$ws.Range("B351").Value = "Huasca De Ocampo"
$ws.Range("B354").Value = "Huejutla De Reyes"
$ws.Range("B357").Value = "Jacala De Ledezma"
$ws.Range("B362").Value = "Mineral Del Chico"
$ws.Range("B363").Value = "Mineral Del Monte"
$ws.Range("B364").Value = "Mixquiahuala De Juárez"
$ws.Range("B365").Value = "Molango De Escamilla"
$ws.Range("B367").Value = "Nopala De Villagrán"
$ws.Range("B368").Value = "Pachuca De Soto"
$ws.Range("B371").Value = "Progreso De Obregón"
$ws.Range("B376").Value = "Santiago De Anaya"
$ws.Range("B380").Value = "Tenango De Doria"
$ws.Range("B382").Value = "Tepehuacán De Guerrero"
$ws.Range("B383").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B386").Value = "Tezontepec De Aldama"
$ws.Range("B393").Value = "Tula De Allende"
$ws.Range("B394").Value = "Tulancingo De Bravo"
$ws.Range("B398").Value = "Zacualtipán De Ángeles"
$ws.Range("B404").Value = "Atotonilco El Alto"
$ws.Range("B405").Value = "Autlán De Navarro"
$ws.Range("B410").Value = "Encarnación De Díaz"
$ws.Range("B412").Value = "Ixtlahuacán Del Río"
$ws.Range("B413").Value = "Jilotlán De Los Dolores"
$ws.Range("B416").Value = "Lagos De Moreno"
$ws.Range("B421").Value = "San Cristóbal De La Barranca"
$ws.Range("B422").Value = "San Juan De Los Lagos"
$ws.Range("B423").Value = "San Miguel El Alto"
$ws.Range("B425").Value = "Tamazula De Gordiano"
$ws.Range("B428").Value = "Teocuitatlán De Corona"
$ws.Range("B429").Value = "Tepatitlán De Morelos"
$ws.Range("B430").Value = "Tizapán El Alto"
$ws.Range("B431").Value = "Tlajomulco De Zúñiga"
$ws.Range("B436").Value = "Unión De San Antonio"
$ws.Range("B437").Value = "Valle De Juárez"
$ws.Range("B440").Value = "Zapotlán El Grande"
$ws.Range("B500").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B518").Value = "Coatlán Del Río"
$ws.Range("B526").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B530").Value = "Puente De Ixtla"
$ws.Range("B536").Value = "Tetela Del Volcán"
$ws.Range("B537").Value = "Tlaltizapán De Zapata"
$ws.Range("B545").Value = "Zacualpan De Amilpas"
$ws.Range("B548").Value = "Amatlán De Cañas"
$ws.Range("B557").Value = "Mier Y Noriega"
$ws.Range("B563").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B568").Value = "Ayoquezco De Aldama"
$ws.Range("B573").Value = "Chalcatongo De Hidalgo"
$ws.Range("B574").Value = "Ciénega De Zimatlán"
$ws.Range("B577").Value = "Coicoyán De Las Flores"
$ws.Range("B580").Value = "Constancia Del Rosario"
$ws.Range("B582").Value = "Cuilápam De Guerrero"
$ws.Range("B583").Value = "Fresnillo De Trujano"
$ws.Range("B584").Value = "Guadalupe De Ramírez"
$ws.Range("B586").Value = "Guevea De Humboldt"
$ws.Range("B587").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B588").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B589").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B591").Value = "Huautla De Jiménez"
$ws.Range("B593").Value = "Ixtlán De Juárez"
$ws.Range("B594").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B602").Value = "Mariscala De Juárez"
$ws.Range("B603").Value = "Mártires De Tacubaya"
$ws.Range("B605").Value = "Mazatlán Villa De Flores"
$ws.Range("B607").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B609").Value = "Nejapa De Madero"
$ws.Range("B610").Value = "Oaxaca De Juárez"
$ws.Range("B611").Value = "Ocotlán De Morelos"
$ws.Range("B612").Value = "Pinotepa De Don Luis"
$ws.Range("B613").Value = "Putla Villa De Guerrero"
$ws.Range("D621").Value = 0.0009689922480620156
$ws.Range("B630").Value = "San Antonino El Alto"
$ws.Range("B632").Value = "San Antonio De La Cal"
$ws.Range("B642").Value = "San Dionisio Del Mar"
$ws.Range("B647").Value = "San Francisco Del Mar"
$ws.Range("B673").Value = "San Juan Del Estado"
$ws.Range("B674").Value = "San Juan Del Río"
$ws.Range("D704").Value = 0.0009302325581395348
$ws.Range("B715").Value = "San Miguel Del Puerto"
$ws.Range("B732").Value = "San Pedro El Alto"
$ws.Range("B759").Value = "Santa Ana Del Valle"
$ws.Range("B768").Value = "Santa Cruz De Bravo"
$ws.Range("B771").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B776").Value = "Santa Inés De Zaragoza"
$ws.Range("B777").Value = "Santa Inés Del Monte"
$ws.Range("B778").Value = "Santa Lucía Del Camino"
$ws.Range("B836").Value = "Santo Domingo De Morelos"
$ws.Range("B848").Value = "Sitio De Xitlapehua"
$ws.Range("B849").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B850").Value = "Tataltepec De Valdés"
$ws.Range("B851").Value = "Teotitlán De Flores Magón"
$ws.Range("B852").Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Range("D852").Value = 0.0009689922480620156
$ws.Range("B853").Value = "Tlacolula De Matamoros"
$ws.Range("B855").Value = "Tlalixtac De Cabrera"
$ws.Range("B856").Value = "Totontepec Villa De Morelos"
$ws.Range("B858").Value = "Villa De Chilapa De Díaz"
$ws.Range("B859").Value = "Villa De Etla"
$ws.Range("B860").Value = "Villa De Tamazulápam Del Progreso"
$ws.Range("B861").Value = "Villa De Tututepec"
$ws.Range("B862").Value = "Villa De Zaachila"
$ws.Range("B864").Value = "Villa Sola De Vega"
$ws.Range("B865").Value = "Zapotitlán Del Río"
$ws.Range("B868").Value = "Zimatlán De Álvarez"
$ws.Range("B895").Value = "Ayotoxco De Guerrero"
$ws.Range("B900").Value = "Chalchicomula De Sesma"
$ws.Range("B910").Value = "Chila De La Sal"
$ws.Range("D917").Value = 0.009534883720930231
$ws.Range("B921").Value = "Cuapiaxtla De Madero"
$ws.Range("B925").Value = "Cuayuca De Andrade"
$ws.Range("B926").Value = "Cuetzalan Del Progreso"
$ws.Range("B941").Value = "Huehuetlán El Chico"
$ws.Range("B942").Value = "Huehuetlán El Grande"
$ws.Range("B946").Value = "Huitzilan De Serdán"
$ws.Range("B948").Value = "Ixcamilpa De Guerrero"
$ws.Range("B952").Value = "Izúcar De Matamoros"
$ws.Range("B962").Value = "Los Reyes De Juárez"
$ws.Range("B963").Value = "Mazapiltepec De Juárez"
$ws.Range("B976").Value = "Palmar De Bravo"
$ws.Range("B986").Value = "San Diego La Mesa Tochimiltzingo"
$ws.Range("B1001").Value = "San Nicolás De Los Ranchos"
$ws.Range("B1005").Value = "San Salvador El Seco"
$ws.Range("B1006").Value = "San Salvador El Verde"
$ws.Range("B1013").Value = "Tecali De Herrera"
$ws.Range("B1021").Value = "Tepanco De López"
$ws.Range("B1022").Value = "Tepango De Rodríguez"
$ws.Range("B1023").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B1029").Value = "Tepexi De Rodríguez"
$ws.Range("B1031").Value = "Tepeyahualco De Cuauhtémoc"
$ws.Range("B1032").Value = "Tetela De Ocampo"
$ws.Range("B1033").Value = "Teteles De Avila Castillo"
$ws.Range("B1038").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B1049").Value = "Totoltepec De Guerrero"
$ws.Range("B1051").Value = "Tuzamapan De Galeana"
$ws.Range("D1053").Value = 0.0009302325581395348
$ws.Range("B1055").Value = "Xayacatlán De Bravo"
$ws.Range("B1061").Value = "Xochitlán De Vicente Suárez"
$ws.Range("D1064").Value = 0.0009689922480620156
$ws.Range("B1077").Value = "Amealco De Bonfil"
$ws.Range("B1078").Value = "Cadereyta De Montes"
$ws.Range("B1081").Value = "Landa De Matamoros"
$ws.Range("B1083").Value = "Pinal De Amoles"
$ws.Range("B1086").Value = "San Juan Del Río"
$ws.Range("B1095").Value = "Ciudad Del Maíz"
$ws.Range("B1100").Value = "Mexquitic De Carmona"
$ws.Range("B1105").Value = "San Ciro De Acosta"
$ws.Range("B1109").Value = "Santa María Del Río"
$ws.Range("B1113").Value = "Tanquián De Escobedo"
$ws.Range("B1116").Value = "Villa De Ramos"
$ws.Range("B1117").Value = "Villa De Reyes"
$ws.Range("D1134").Value = 0.0009689922480620156
$ws.Range("B1173").Value = "Acuamanala De Miguel Hidalgo"
$ws.Range("B1175").Value = "Amaxac De Guerrero"
$ws.Range("B1180").Value = "Contla De Juan Cuamatzi"
$ws.Range("D1182").Value = 0.0009689922480620156
$ws.Range("B1184").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B1188").Value = "Muñoz De Domingo Arenas"
$ws.Range("B1189").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B1192").Value = "Papalotla De Xicohténcatl"
$ws.Range("B1197").Value = "San Pablo Del Monte"
$ws.Range("D1205").Value = 0.0009689922480620156
$ws.Range("B1206").Value = "Tepetitla De Lardizábal"
$ws.Range("B1209").Value = "Tetla De La Solidaridad"
$ws.Range("B1219").Value = "Ziltlaltépec De Trinidad Sánchez Santos"
$ws.Range("B1227").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B1231").Value = "Amatlán De Los Reyes"
$ws.Range("B1238").Value = "Boca Del Río"
$ws.Range("B1242").Value = "Castillo De Teayo"
$ws.Range("B1244").Value = "Cazones De Herrera"
$ws.Range("D1253").Value = 0.0009689922480620156
$ws.Range("B1257").Value = "Cosamaloapan De Carpio"
$ws.Range("B1272").Value = "Hueyapan De Ocampo"
$ws.Range("B1273").Value = "Huiloapan De Cuauhtémoc"
$ws.Range("B1277").Value = "Ixhuatlán De Madero"
$ws.Range("B1278").Value = "Ixhuatlán Del Café"
$ws.Range("B1286").Value = "Juchique De Ferrer"
$ws.Range("B1290").Value = "Las Vigas De Ramírez"
$ws.Range("B1291").Value = "Lerdo De Tejada"
$ws.Range("B1294").Value = "Martínez De La Torre"
$ws.Range("B1295").Value = "Medellín De Bravo"
$ws.Range("B1298").Value = "Mixtla De Altamirano"
$ws.Range("B1309").Value = "Paso De Ovejas"
$ws.Range("B1310").Value = "Paso Del Macho"
$ws.Range("B1314").Value = "Poza Rica De Hidalgo"
$ws.Range("B1323").Value = "Sayula De Alemán"
$ws.Range("B1325").Value = "Soledad De Doblado"
$ws.Range("B1352").Value = "Vega De Alatorre"
$ws.Range("B1361").Value = "Zontecomatlán De López Y Fuentes"
$ws.Range("B1362").Value = "Zozocolco De Hidalgo"
$ws.Range("B1373").Value = "Jiménez Del Teul"
$ws.Range("B1375").Value = "Nochistlán De Mejía"
$ws.Range("B1381").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B1383").Value = "Villa De Cos"

$ws.Rows("1388:1392").Delete()

$wb.Save()